# Generate Report for Handoff
#
# The e2e\ed85a594-4fb9-4fe9-8d41-36aef1ab5cdc.md file has just been handed
# off again, so refresh its "Latest Handoff / Generate" timestamps on each
# sheet:
#   - Overview!G7            -> Latest HO Xliff Generate Date
#   - zh-cn!H7 (row 7 data)  -> Latest Handoff Datetime
#   - de-de!H7 (row 7 data)  -> Latest Handoff Datetime

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-30 12:51:40"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-30 12:51:36"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-30 12:51:40"
